# Adds the new "2024" worksheet with refreshed PBO cashflow data and updates the
# "2023" worksheet's selection, matching the 12/30/2023 LDI data refresh.

$wb = $excel.ActiveWorkbook

# The "2023" sheet is currently the last sheet and is used as a template for the new
# "2024" sheet so that column styles (e.g. the date format on column A), the column
# width, and the column E "=SUM(Bn:Dn)" formulas are preserved exactly.
$ws2023 = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws2023.Copy([System.Reflection.Missing]::Value, $ws2023) | Out-Null
$ws2024 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2024.Name = "2024"

# New PBO cashflow data for 2024 - columns A (date), B (Retirement), C (Pension),
# D (IBT). Column E keeps the inherited "=SUM(Bn:Dn)" formulas copied from the 2023
# sheet, so Excel recalculates the totals from the new inputs automatically.
$jagged2024 = @(
    @(45291,72771.80032437398,177464.42698958149,965775.58824944135),
    @(45657,306159.98416596366,750911.73545306374,3528912.8894307534),
    @(46022,671576.52076659957,1608494.8015452409,6739523.1389516024),
    @(46387,1066998.7043058577,2577009.1786195692,9935049.5052611958),
    @(46752,1483701.3053002157,3653608.475762838,13097052.150619233),
    @(47118,1936560.7234457438,4763453.2448138297,16201407.37712601),
    @(47483,2375935.6398847671,5834135.6156361299,19204572.72648403),
    @(47848,2757369.9773275829,6894007.0482303184,22132462.933259532),
    @(48213,3091861.4745630282,7947923.4681398747,25003215.70856474),
    @(48579,3423925.1694616345,8980680.332138367,27824584.971314386),
    @(48944,3729805.2970696134,9971165.054259561,30580692.060145281),
    @(49309,4003182.0646350044,10920236.163193503,33246237.106070064),
    @(49674,4275328.1865917668,11836946.42913226,35837952.596104912),
    @(50040,4527364.3223877735,12712192.178676147,38378570.324450068),
    @(50405,4764604.178078956,13555630.122654112,40867671.32649491),
    @(50770,4989474.6013234546,14379645.789431527,43304573.323545307),
    @(51135,5182859.7178554283,15185220.575955097,45695836.945055395),
    @(51501,5362464.2344171181,15967586.541546524,48034798.455710277),
    @(51866,5528674.4011129607,16733273.470331112,50288773.511698306),
    @(52231,5680719.8152656127,17495802.950941537,52441919.461787142),
    @(52596,5807227.3568064841,18246113.593777481,54510363.164517142),
    @(52962,5904517.6440829635,18940054.084129229,56479798.40966434),
    @(53327,5988781.5076424368,19573697.422852647,58344190.777100928),
    @(53692,6045074.3088858854,20176039.311403371,60100677.154142365),
    @(54057,6069364.3391083954,20753035.969552781,61707511.704928845),
    @(54423,6066575.9212257564,21295627.636774462,63149062.290190779),
    @(54788,6029918.3699446227,21774787.24836291,64423260.792275697),
    @(55153,5961628.1602434553,22180954.787842605,65508674.302369557),
    @(55518,5850916.2273299815,22520034.080789819,66373276.454252981),
    @(55884,5694231.5568855843,22821907.872636441,66991639.169685028),
    @(56249,5506346.2136883875,23105652.342151131,67348439.949188411),
    @(56614,5303254.8292722097,23342592.337297857,67413829.020718157),
    @(56979,5085208.3644563155,23527732.891550135,67187111.68767333),
    @(57345,4850206.4343097173,23686075.303323604,66662668.354720213),
    @(57710,4604346.8834955916,23859625.947518785,65813220.407746866),
    @(58075,4351728.3801072128,24084389.331684839,64652362.6961568),
    @(58440,4097303.8126994427,24345718.374628063,63199946.878036387),
    @(58806,3842104.5548687773,24639966.779594369,61489155.945694536),
    @(59171,3586865.1671211724,24926319.818440434,59574635.611837476),
    @(59536,3334570.7960938327,24950781.12177293,57497147.830955565),
    @(59901,3087811.9407089963,24727321.375946153,55288778.940996267),
    @(60267,2847926.4608056676,24541530.244220577,52980695.287904076),
    @(60632,2616005.9155070838,24471479.741723377,50600842.566382915),
    @(60997,2392916.8145932881,24144067.705091767,48180283.60102462),
    @(61362,2179314.0940557346,23347516.365740295,45745818.437779978),
    @(61728,1975685.8753442138,22399021.474281408,43314980.064678378),
    @(62093,1782369.4618177321,21449904.030640997,40899138.515155621),
    @(62458,1599610.7079457934,20501508.229541045,38511337.550541446),
    @(62823,1427577.6908498642,19554646.115213819,36159984.615720049),
    @(63189,1266391.6568639553,18610733.731921304,33852582.616286203),
    @(63554,1116143.6790025788,17672103.583818275,31598306.584025946),
    @(63919,976894.32707455894,16740221.484414572,29403005.907350045),
    @(64284,848673.34711659106,15816233.107635615,27269560.982332487),
    @(64650,731449.83350302675,14901327.722096797,25200939.369909637),
    @(65015,625117.50209478172,13996770.757690893,23199989.018412698),
    @(65380,529481.65742549626,13103931.270369213,21269471.504213531),
    @(65745,444250.77907996526,12224350.481825596,19412442.149142299),
    @(66111,369045.20742421946,11359762.007459346,17632181.722071242),
    @(66476,303375.62693369086,10512050.862408891,15932051.471404547),
    @(66841,246672.57897391249,9683224.8919474725,14315480.741325365),
    @(67206,198298.33444524027,8875389.4219531976,12785749.822815048),
    @(67572,157551.99018909762,8090753.1432746788,11345826.222618757),
    @(67937,123683.14037814221,7331588.994704783,9998283.3281803057),
    @(68302,95908.424161760558,6600263.3273885632,8745237.5196869206),
    @(68667,73440.373367095483,5899327.2855516253,7588221.5423531858),
    @(69033,55518.091189105384,5231517.7382389288,6528136.6722141728),
    @(69398,41425.481155556015,4599683.4440315571,5565177.4500213917),
    @(69763,30503.164132491031,4006838.9937306931,4698621.5718206642),
    @(70128,22157.705029166562,3456032.8391464073,3926671.2394719822),
    @(70494,15872.660199364222,2949927.0439884998,3246390.9214468687),
    @(70859,11208.533871991396,2490358.9631513627,2653649.3664118415),
    @(71224,7798.4210700492404,2078138.5883057974,2143197.3756314297),
    @(71589,5343.0363239127346,1712995.086175415,1708961.7942459458),
    @(71955,3602.4726471449876,1393522.2886637126,1344288.1993869469),
    @(72320,2388.2343735553541,1117476.7053035055,1042166.2714711443),
    @(72685,1555.4265660271103,882293.87538093468,795513.46712459926),
    @(73050,994.26323043094726,685089.4248500024,597326.51804733765),
    @(73415,623.11627787904831,522537.33866159775,440755.36914216698),
    @(73780,382.50768883598539,391005.69133188727,319287.96559253283),
    @(74145,229.76988138846303,286676.53228185087,226858.86975870773)
)

$dataRows = $jagged2024.Count
$dataCols = $jagged2024[0].Count
$arr2024 = New-Object 'object[,]' $dataRows,$dataCols
for ($r = 0; $r -lt $dataRows; $r++) {
    for ($c = 0; $c -lt $dataCols; $c++) {
        $arr2024[$r,$c] = $jagged2024[$r][$c]
    }
}

$ws2024.Range("A2:D81").Value = $arr2024

# Refresh the selection/active-cell state to match the saved view.
$ws2024.Range("I11").Select() | Out-Null

$ws2023.Activate() | Out-Null
$ws2023.Range("C40").Select() | Out-Null

$ws2024.Activate() | Out-Null
